$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "t5tss4s"
$ws.Range("B2").Value = "agAmuvA"

$ws.Range("A3").Value = "mngr198435"
$ws.Range("B3").Value = "ih7677r77"

$ws.Range("A4").Value = "rtdtdttd"
$ws.Range("B4").Value = "frydyydydy"

$ws.Range("B4").Select()
